$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 366.85715
$ws.Range("I115").Value = 366.85715
$ws.Range("K115").Value = 1100.57145
$ws.Range("M115").Value = 466.4285500000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1063.8
$ws.Range("I97").Value = 981.0714
$ws.Range("J97").Value = 2222
$ws.Range("K97").Value = 981.0714
$ws.Range("L97").Value = 2222
$ws.Range("M97").Value = -485.0714
$ws.Range("N97").Value = -3214
$ws.Range("H102").Value = 1851.8
$ws.Range("I102").Value = 1827.5
$ws.Range("K102").Value = 1827.5
$ws.Range("M102").Value = -205.5
$ws.Range("H122").Value = 5973.0835
$ws.Range("I122").Value = 6299.1816
$ws.Range("K122").Value = 18897.5448
$ws.Range("M122").Value = -16447.5448
$ws.Range("H132").Value = 2952.5881
$ws.Range("I132").Value = 3181.4443
$ws.Range("K132").Value = 9544.332900000001
$ws.Range("M132").Value = -7014.332900000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 517.4167
$ws.Range("I7").Value = 572.35297
$ws.Range("J7").Value = 384
$ws.Range("K7").Value = 572.35297
$ws.Range("L7").Value = 384
$ws.Range("M7").Value = -459.35297
$ws.Range("N7").Value = -610
$ws.Range("H16").Value = 8104.143
$ws.Range("I16").Value = 3025.875
$ws.Range("K16").Value = 3025.875
$ws.Range("M16").Value = -2738.875
$ws.Range("H31").Value = 39377.75
$ws.Range("I31").Value = 2253.111
$ws.Range("J31").Value = 56963.105
$ws.Range("K31").Value = 2253.111
$ws.Range("L31").Value = 56963.105
$ws.Range("M31").Value = -1958.111
$ws.Range("N31").Value = -57553.105
$ws.Range("H34").Value = 39377.75
$ws.Range("I34").Value = 2253.111
$ws.Range("J34").Value = 56963.105
$ws.Range("K34").Value = 2253.111
$ws.Range("L34").Value = 56963.105
$ws.Range("M34").Value = -2051.111
$ws.Range("N34").Value = -57367.105
$ws.Range("H62").Value = 3861
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 8000
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 3861
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 40000
$ws.Range("N65").Value = -46240
$ws.Range("H94").Value = 2975
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H105").Value = 542
$ws.Range("I105").Value = 560.64703
$ws.Range("K105").Value = 560.64703
$ws.Range("M105").Value = 1186.35297
$ws.Range("H113").Value = 8104.143
$ws.Range("I113").Value = 3025.875
$ws.Range("K113").Value = 3025.875
$ws.Range("M113").Value = -855.875
$ws.Range("H134").Value = 229463.89
$ws.Range("I134").Value = 2242.1162
$ws.Range("K134").Value = 6726.348599999999
$ws.Range("M134").Value = -4191.348599999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1123.091
$ws.Range("I134").Value = 1123.091
$ws.Range("K134").Value = 3369.273
$ws.Range("M134").Value = 1700.727
$ws.Range("H139").Value = 7214.8076
$ws.Range("I139").Value = 6640.933
$ws.Range("K139").Value = 19922.799
$ws.Range("M139").Value = -14782.799

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 58430
$ws.Range("J32").Value = 58430
$ws.Range("L32").Value = 58430
$ws.Range("N32").Value = -59022
$ws.Range("H100").Value = 54855
$ws.Range("J100").Value = 54855
$ws.Range("L100").Value = 54855
$ws.Range("N100").Value = -57019
$ws.Range("H113").Value = 508112.5
$ws.Range("I113").Value = 1430485.6
$ws.Range("J113").Value = 11450.077
$ws.Range("K113").Value = 1430485.6
$ws.Range("L113").Value = 11450.077
$ws.Range("M113").Value = -1428315.6
$ws.Range("N113").Value = -15790.077
$ws.Range("H122").Value = 4786.25
$ws.Range("I122").Value = 2695
$ws.Range("J122").Value = 5483.3335
$ws.Range("K122").Value = 8085
$ws.Range("L122").Value = 16450.0005
$ws.Range("M122").Value = -5635
$ws.Range("N122").Value = -21350.0005
$ws.Range("H123").Value = 47124.375
$ws.Range("J123").Value = 47124.375
$ws.Range("L123").Value = 47124.375
$ws.Range("N123").Value = -52024.375
$ws.Range("H131").Value = 30326
$ws.Range("J131").Value = 30326
$ws.Range("L131").Value = 30326
$ws.Range("N131").Value = -40406

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4616
$ws.Range("I93").Value = 4499.5
$ws.Range("K93").Value = 4499.5
$ws.Range("M93").Value = -3251.5
$ws.Range("H122").Value = 6667169.5
$ws.Range("I122").Value = 5000752
$ws.Range("J122").Value = 10000005
$ws.Range("K122").Value = 15002256
$ws.Range("L122").Value = 30000015
$ws.Range("M122").Value = -14999806
$ws.Range("N122").Value = -30004915

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 2249.5
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H31").Value = 3000
$ws.Range("J31").Value = 3000
$ws.Range("L31").Value = 3000
$ws.Range("N31").Value = -3696
$ws.Range("H62").Value = 83769.38
$ws.Range("I62").Value = 257249.75
$ws.Range("J62").Value = 6667
$ws.Range("K62").Value = 257249.75
$ws.Range("L62").Value = 6667
$ws.Range("M62").Value = -256625.75
$ws.Range("N62").Value = -7915
$ws.Range("H65").Value = 83769.38
$ws.Range("I65").Value = 257249.75
$ws.Range("J65").Value = 6667
$ws.Range("K65").Value = 1286248.75
$ws.Range("L65").Value = 33335
$ws.Range("M65").Value = -1283128.75
$ws.Range("N65").Value = -39575
$ws.Range("H81").Value = 2065.077
$ws.Range("I81").Value = 1584.4
$ws.Range("J81").Value = 3667.3333
$ws.Range("K81").Value = 3168.8
$ws.Range("L81").Value = 7334.6666
$ws.Range("M81").Value = -2107.8
$ws.Range("N81").Value = -9456.6666
$ws.Range("H84").Value = 2065.077
$ws.Range("I84").Value = 1584.4
$ws.Range("J84").Value = 3667.3333
$ws.Range("K84").Value = 15844
$ws.Range("L84").Value = 36673.333
$ws.Range("M84").Value = -10540
$ws.Range("N84").Value = -47281.333
$ws.Range("H126").Value = 1740.4
$ws.Range("I126").Value = 1740.4
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5221.200000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2751.200000000001
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 10467035
$ws.Range("I136").Value = 13736026
$ws.Range("K136").Value = 41208078
$ws.Range("M136").Value = -41205528
